# Fixed Bento 80 Test scripts
# Appends "order By ... LIMIT 100" clauses to the three Cypher queries stored
# in column B (the "query" column) of the "startup" sheet, matching the
# author's fix to the TC13_Bento_Filter_Chemo-TAC.xlsx generator scripts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("CasesTab" row) - query that actually returns per-case fields;
# append an ORDER BY / LIMIT clause after the `Survival (days)` column.
$b2 = $ws.Cells.Item(2, 2).Value2
$ws.Cells.Item(2, 2).Value = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100 "

# Row 3 ("SamplesTab" row) - query that returns per-sample fields; append an
# ORDER BY / LIMIT clause after the `Sample Procurement Method` column.
$b3 = $ws.Cells.Item(3, 2).Value2
$ws.Cells.Item(3, 2).Value = $b3 + "`n order By samp.sample_id ASC LIMIT 100"

# Row 4 ("FilesTab" row) - query that lists files; replace the trailing
# "    order by f.file_name" with the new capitalized clause plus LIMIT.
$b4 = $ws.Cells.Item(4, 2).Value2
$b4 = $b4.Replace("    order by f.file_name", "     order By f.file_name ASC LIMIT 100")
$ws.Cells.Item(4, 2).Value = $b4

# Update row heights to account for the extra line of text that was added
# to rows 2 and 3 (row 4 was already at Excel's row-height cap).
$ws.Rows.Item(2).RowHeight = 331.2
$ws.Rows.Item(3).RowHeight = 360

# Match the author's final selection/scroll state: B4 was the last cell
# edited, so it ends up as the active cell/selection.
$ws.Range("B4").Select()
